$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 30 (CUGerm): fix citation typo "seeds:_1994" -> "seeds_1994"
$ws.Range("F30").Value = '[@heydecker_seed_1972; @bewley_seeds_1994]'

# Row 13 (CVG / Kotowski's coefficient of velocity): fix citation typo "seeds:_1994" -> "seeds_1994"
$ws.Range("F13").Value = '[@kotowski_temperature_1926, @nichols_two_1968; @bewley_seeds_1994; @labouriau_uma_1983; @scott_review_1984]'

# Row 23 (TimsonsIndex): fix citation typo "seeds:_1998" -> "seeds_1998"
$ws.Range("F23").Value = '[@grose_laboratory_1958; @timson_new_1965; @brown_representing_1988; @baskin_seeds_1998; @goodchild_method_1971]'

# Row 2 (GermPercent): expand the Germination index description to include "Final germination percentage or"
$ws.Range("A2").Value = 'Germination percentage or Final germination percentage or Germinability ($GP$)'

# Reflect the final saved view: active cell A2, no pinned top-left row
$ws.Activate()
$ws.Range("A2").Select()
